$d = $word.ActiveDocument

# Update the heading date paragraph
$d.Paragraphs.Item(1).Range.Text = "2023-05-02 Tuesday"

# Update each multiplication-problem cell in document order (row-major, 5 cols/row)
$t = $d.Tables.Item(1)
$values = @(
  "13×28=",
  "42×78=",
  "18×90=",
  "54×43=",
  "16×44=",
  "83×45=",
  "31×28=",
  "35×15=",
  "23×18=",
  "66×24=",
  "71×87=",
  "25×26=",
  "46×57=",
  "40×12=",
  "93×20=",
  "99×10=",
  "45×32=",
  "79×63=",
  "100×49=",
  "90×57=",
  "58×12=",
  "31×48=",
  "13×35=",
  "31×51=",
  "25×86=",
  "12×31=",
  "46×15=",
  "16×92=",
  "59×59=",
  "57×93=",
  "76×59=",
  "83×17=",
  "16×46=",
  "57×94=",
  "84×42=",
  "63×15=",
  "59×46=",
  "44×47=",
  "22×42=",
  "29×64=",
  "57×36=",
  "69×27=",
  "77×51=",
  "43×12=",
  "70×36=",
  "50×98=",
  "96×21=",
  "96×93=",
  "52×38=",
  "72×46=",
  "74×73=",
  "88×38=",
  "49×11=",
  "16×24=",
  "99×38=",
  "17×99=",
  "68×57=",
  "16×15=",
  "14×14=",
  "54×89=",
  "50×48=",
  "66×18=",
  "21×99=",
  "85×72=",
  "17×15=",
  "51×39=",
  "28×56=",
  "74×81=",
  "73×55=",
  "66×31=",
  "68×45=",
  "65×47=",
  "67×47=",
  "66×99=",
  "82×72=",
  "94×49=",
  "86×27=",
  "80×97=",
  "13×56=",
  "47×63=",
  "76×77=",
  "10×82=",
  "80×64=",
  "95×89=",
  "55×87=",
  "78×18=",
  "67×23=",
  "50×69=",
  "98×31=",
  "87×10=",
  "20×20=",
  "63×37=",
  "90×95=",
  "65×57=",
  "14×90=",
  "15×36=",
  "75×40=",
  "99×23=",
  "28×68=",
  "32×35="
)

$cols = $t.Columns.Count
for ($i = 0; $i -lt $values.Count; $i++) {
  $row = [math]::Floor($i / $cols) + 1
  $col = ($i % $cols) + 1
  $t.Cell($row, $col).Range.Text = $values[$i]
}

Write-Host "Done updating" $values.Count "cells"
